$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.665.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.911.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.00%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.553'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.908.54'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.86'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.87%  '
$ws.Range('E12').Value = '  -4.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.413.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.776.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.905.99'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.99%  '
$ws.Range('E18').Value = '  -7.42%  '
$ws.Range('E19').Value = '  -3.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.53'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '360.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.036.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.450'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.181'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0861'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -13.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -11.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.67'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.68'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.35'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.28%  '
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.34%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.335.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.642'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.81'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0569'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.83'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.34%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0233'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.00%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0923'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.95%  '
